$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / header text updates
$ws.Range("A2").Value = "Year-to-Date through November 2016 and November 2015 (Thousand Tons)"

# Column header row 5: October -> November YTD labels
$ws.Range("B5").Value = "November 2016 YTD"
$ws.Range("E5").Value = "November 2016 YTD"
$ws.Range("G5").Value = "November 2016 YTD"
$ws.Range("I5").Value = "November 2016 YTD"
$ws.Range("K5").Value = "November 2016 YTD"
$ws.Range("C5").Value = "November 2015 YTD"
$ws.Range("F5").Value = "November 2015 YTD"
$ws.Range("H5").Value = "November 2015 YTD"
$ws.Range("J5").Value = "November 2015 YTD"
$ws.Range("L5").Value = "November 2015 YTD"

# Data cell updates (rows 6-67)
$ws.Range("B6").Value = 973
$ws.Range("C6").Value = 1709
$ws.Range("D6").Value = -0.43
$ws.Range("E6").Value = 124
$ws.Range("F6").Value = 398
$ws.Range("G6").Value = 843
$ws.Range("H6").Value = 1300
$ws.Range("K6").Value = 7

$ws.Range("B7").Value = 74
$ws.Range("D7").Value = -0.79
$ws.Range("G7").Value = 74

$ws.Range("B8").Value = 14
$ws.Range("C8").Value = 20
$ws.Range("D8").Value = -0.3
$ws.Range("G8").Value = 11
$ws.Range("L8").Value = 8

$ws.Range("B9").Value = 762
$ws.Range("C9").Value = 932
$ws.Range("D9").Value = -0.18
$ws.Range("G9").Value = 759
$ws.Range("H9").Value = 929

$ws.Range("B10").Value = 124
$ws.Range("C10").Value = 398
$ws.Range("E10").Value = 124
$ws.Range("F10").Value = 398

$ws.Range("B13").Value = 24351
$ws.Range("C13").Value = 30697
$ws.Range("G13").Value = 24158
$ws.Range("H13").Value = 30478
$ws.Range("K13").Value = 193
$ws.Range("L13").Value = 218

$ws.Range("B14").Value = 504
$ws.Range("C14").Value = 711
$ws.Range("D14").Value = -0.29
$ws.Range("G14").Value = 504
$ws.Range("H14").Value = 711

$ws.Range("B15").Value = 691
$ws.Range("C15").Value = 1056
$ws.Range("D15").Value = -0.35
$ws.Range("H15").Value = 995
$ws.Range("K15").Value = 64
$ws.Range("L15").Value = 61

$ws.Range("B16").Value = 23156
$ws.Range("C16").Value = 28931
$ws.Range("G16").Value = 23027
$ws.Range("H16").Value = 28773
$ws.Range("K16").Value = 129
$ws.Range("L16").Value = 157

$ws.Range("B17").Value = 130080
$ws.Range("C17").Value = 154291
$ws.Range("E17").Value = 78820
$ws.Range("F17").Value = 90009
$ws.Range("G17").Value = 50525
$ws.Range("H17").Value = 63448
$ws.Range("I17").Value = 25
$ws.Range("J17").Value = 36
$ws.Range("K17").Value = 709
$ws.Range("L17").Value = 797

$ws.Range("B18").Value = 32549
$ws.Range("C18").Value = 41385
$ws.Range("E18").Value = 1759
$ws.Range("F18").Value = 1882
$ws.Range("G18").Value = 30278
$ws.Range("H18").Value = 38971
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 12
$ws.Range("K18").Value = 501
$ws.Range("L18").Value = 521

$ws.Range("B19").Value = 32578
$ws.Range("C19").Value = 36669
$ws.Range("D19").Value = -0.11
$ws.Range("E19").Value = 30930
$ws.Range("F19").Value = 34118
$ws.Range("G19").Value = 1636
$ws.Range("H19").Value = 2542
$ws.Range("I19").Value = 11
$ws.Range("J19").Value = 9
$ws.Range("L19").Value = 0.19

$ws.Range("B20").Value = 21160
$ws.Range("C20").Value = 27321
$ws.Range("E20").Value = 20871
$ws.Range("F20").Value = 27032
$ws.Range("G20").Value = 228
$ws.Range("H20").Value = 209
$ws.Range("K20").Value = 57
$ws.Range("L20").Value = 65

$ws.Range("B21").Value = 26310
$ws.Range("C21").Value = 28919
$ws.Range("D21").Value = -0.09
$ws.Range("E21").Value = 7887
$ws.Range("F21").Value = 7140
$ws.Range("G21").Value = 18383
$ws.Range("H21").Value = 21726
$ws.Range("K21").Value = 40
$ws.Range("L21").Value = 53

$ws.Range("B22").Value = 17482
$ws.Range("C22").Value = 19996
$ws.Range("D22").Value = -0.13
$ws.Range("E22").Value = 17373
$ws.Range("F22").Value = 19838
$ws.Range("K22").Value = 109
$ws.Range("L22").Value = 158

$ws.Range("B23").Value = 104888
$ws.Range("C23").Value = 115636
$ws.Range("D23").Value = -0.093
$ws.Range("E23").Value = 103663
$ws.Range("F23").Value = 114146
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 12
$ws.Range("I23").Value = 52
$ws.Range("J23").Value = 51
$ws.Range("K23").Value = 1164
$ws.Range("L23").Value = 1427

$ws.Range("B24").Value = 13537
$ws.Range("C24").Value = 16520
$ws.Range("D24").Value = -0.18
$ws.Range("E24").Value = 12885
$ws.Range("F24").Value = 15764
$ws.Range("I24").Value = 33
$ws.Range("J24").Value = 30
$ws.Range("K24").Value = 620
$ws.Range("L24").Value = 726

$ws.Range("B25").Value = 13163
$ws.Range("C25").Value = 14862
$ws.Range("D25").Value = -0.11
$ws.Range("E25").Value = 13163
$ws.Range("F25").Value = 14862

$ws.Range("B26").Value = 12520
$ws.Range("C26").Value = 13428
$ws.Range("D26").Value = -0.068
$ws.Range("E26").Value = 12290
$ws.Range("F26").Value = 13121
$ws.Range("I26").Value = 6
$ws.Range("K26").Value = 224
$ws.Range("L26").Value = 302

$ws.Range("B27").Value = 32588
$ws.Range("C27").Value = 35541
$ws.Range("D27").Value = -0.083
$ws.Range("E27").Value = 32552
$ws.Range("F27").Value = 35500
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = 12
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = 13

$ws.Range("B28").Value = 12031
$ws.Range("C28").Value = 13646
$ws.Range("D28").Value = -0.12
$ws.Range("E28").Value = 11780
$ws.Range("F28").Value = 13325
$ws.Range("K28").Value = 251
$ws.Range("L28").Value = 320

$ws.Range("B29").Value = 19797
$ws.Range("C29").Value = 20782
$ws.Range("E29").Value = 19741
$ws.Range("F29").Value = 20716
$ws.Range("K29").Value = 56
$ws.Range("L29").Value = 67

$ws.Range("B30").Value = 1252
$ws.Range("C30").Value = 858
$ws.Range("D30").Value = 0.46
$ws.Range("E30").Value = 1252
$ws.Range("F30").Value = 858

$ws.Range("B31").Value = 94370
$ws.Range("C31").Value = 99272
$ws.Range("D31").Value = -0.049
$ws.Range("E31").Value = 82725
$ws.Range("F31").Value = 87239
$ws.Range("G31").Value = 11380
$ws.Range("H31").Value = 11649
$ws.Range("I31").Value = 17
$ws.Range("J31").Value = 21
$ws.Range("K31").Value = 247
$ws.Range("L31").Value = 363

$ws.Range("C32").Value = 274
$ws.Range("D32").Value = -0.2
$ws.Range("H32").Value = 274

$ws.Range("B34").Value = 16120
$ws.Range("C34").Value = 17849
$ws.Range("D34").Value = -0.097
$ws.Range("E34").Value = 15813
$ws.Range("F34").Value = 17394
$ws.Range("G34").Value = 274
$ws.Range("H34").Value = 423
$ws.Range("K34").Value = 33
$ws.Range("L34").Value = 32

$ws.Range("B35").Value = 17784
$ws.Range("C35").Value = 18455
$ws.Range("D35").Value = -0.036
$ws.Range("E35").Value = 17740
$ws.Range("F35").Value = 18410
$ws.Range("K35").Value = 44
$ws.Range("L35").Value = 45

$ws.Range("B36").Value = 5511
$ws.Range("C36").Value = 5820
$ws.Range("D36").Value = -0.053
$ws.Range("G36").Value = 5492
$ws.Range("H36").Value = 5791
$ws.Range("K36").Value = 19
$ws.Range("L36").Value = 29

$ws.Range("B37").Value = 13618
$ws.Range("C37").Value = 15063
$ws.Range("D37").Value = -0.096
$ws.Range("E37").Value = 13453
$ws.Range("F37").Value = 14874
$ws.Range("G37").Value = 123
$ws.Range("H37").Value = 145
$ws.Range("I37").Value = 10
$ws.Range("J37").Value = 13
$ws.Range("K37").Value = 31
$ws.Range("L37").Value = 30

$ws.Range("B38").Value = 7876
$ws.Range("C38").Value = 8771
$ws.Range("E38").Value = 7836
$ws.Range("F38").Value = 8721
$ws.Range("K38").Value = 40
$ws.Range("L38").Value = 50

$ws.Range("B39").Value = 6758
$ws.Range("C39").Value = 7053
$ws.Range("D39").Value = -0.042
$ws.Range("E39").Value = 6434
$ws.Range("F39").Value = 6605
$ws.Range("G39").Value = 263
$ws.Range("H39").Value = 385
$ws.Range("K39").Value = 53
$ws.Range("L39").Value = 56

$ws.Range("B40").Value = 26482
$ws.Range("C40").Value = 25986
$ws.Range("E40").Value = 21449
$ws.Range("F40").Value = 21236
$ws.Range("G40").Value = 5007
$ws.Range("H40").Value = 4630
$ws.Range("L40").Value = 121

$ws.Range("B41").Value = 63512
$ws.Range("C41").Value = 70495
$ws.Range("D41").Value = -0.099
$ws.Range("E41").Value = 60522
$ws.Range("F41").Value = 67293
$ws.Range("G41").Value = 2807
$ws.Range("H41").Value = 2995
$ws.Range("K41").Value = 182
$ws.Range("L41").Value = 207

$ws.Range("B42").Value = 16008
$ws.Range("C42").Value = 19725
$ws.Range("D42").Value = -0.19
$ws.Range("E42").Value = 15990
$ws.Range("F42").Value = 19703
$ws.Range("K42").Value = 18
$ws.Range("L42").Value = 22

$ws.Range("B43").Value = 29199
$ws.Range("C43").Value = 31926
$ws.Range("D43").Value = -0.085
$ws.Range("E43").Value = 29199
$ws.Range("F43").Value = 31926

$ws.Range("B44").Value = 4159
$ws.Range("C44").Value = 4672
$ws.Range("D44").Value = -0.11
$ws.Range("E44").Value = 1352
$ws.Range("F44").Value = 1677
$ws.Range("G44").Value = 2807
$ws.Range("H44").Value = 2995

$ws.Range("B45").Value = 14145
$ws.Range("C45").Value = 14172
$ws.Range("D45").Value = -0.002
$ws.Range("E45").Value = 13981
$ws.Range("F45").Value = 13988
$ws.Range("K45").Value = 164
$ws.Range("L45").Value = 185

$ws.Range("B46").Value = 108272
$ws.Range("C46").Value = 118335
$ws.Range("D46").Value = -0.085
$ws.Range("E46").Value = 51905
$ws.Range("F46").Value = 58397
$ws.Range("G46").Value = 56209
$ws.Range("H46").Value = 59758
$ws.Range("K46").Value = 159
$ws.Range("L46").Value = 180

$ws.Range("B47").Value = 12311
$ws.Range("C47").Value = 12180
$ws.Range("D47").Value = 0.011
$ws.Range("E47").Value = 10045
$ws.Range("F47").Value = 10054
$ws.Range("G47").Value = 2256
$ws.Range("H47").Value = 2114
$ws.Range("K47").Value = 11
$ws.Range("L47").Value = 13

$ws.Range("B48").Value = 7665
$ws.Range("C48").Value = 10081
$ws.Range("D48").Value = -0.24
$ws.Range("E48").Value = 5242
$ws.Range("F48").Value = 6135
$ws.Range("G48").Value = 2423
$ws.Range("H48").Value = 3945

$ws.Range("B49").Value = 10932
$ws.Range("C49").Value = 14964
$ws.Range("D49").Value = -0.27
$ws.Range("E49").Value = 9693
$ws.Range("F49").Value = 13698
$ws.Range("G49").Value = 1091
$ws.Range("H49").Value = 1098
$ws.Range("K49").Value = 148
$ws.Range("L49").Value = 168

$ws.Range("B50").Value = 77363
$ws.Range("C50").Value = 81110
$ws.Range("D50").Value = -0.046
$ws.Range("E50").Value = 26925
$ws.Range("F50").Value = 28509
$ws.Range("G50").Value = 50439
$ws.Range("H50").Value = 52601

$ws.Range("B51").Value = 81992
$ws.Range("C51").Value = 93688
$ws.Range("D51").Value = -0.12
$ws.Range("E51").Value = 72082
$ws.Range("F51").Value = 82721
$ws.Range("G51").Value = 9526
$ws.Range("H51").Value = 10547
$ws.Range("K51").Value = 383
$ws.Range("L51").Value = 420

$ws.Range("B52").Value = 15025
$ws.Range("C52").Value = 18318
$ws.Range("E52").Value = 15025
$ws.Range("F52").Value = 18318

$ws.Range("B53").Value = 14994
$ws.Range("C53").Value = 15883
$ws.Range("D53").Value = -0.056
$ws.Range("E53").Value = 14978
$ws.Range("F53").Value = 15864
$ws.Range("G53").Value = 13
$ws.Range("H53").Value = 14
$ws.Range("L53").Value = 5

$ws.Range("B54").Value = 12
$ws.Range("C54").Value = 15
$ws.Range("D54").Value = -0.18
$ws.Range("K54").Value = 12
$ws.Range("L54").Value = 15

$ws.Range("B55").Value = 8550
$ws.Range("C55").Value = 9399
$ws.Range("D55").Value = -0.09
$ws.Range("E55").Value = 183
$ws.Range("F55").Value = 195
$ws.Range("G55").Value = 8360
$ws.Range("H55").Value = 9199
$ws.Range("K55").Value = 7
$ws.Range("L55").Value = 5

$ws.Range("B56").Value = 1144
$ws.Range("C56").Value = 1404
$ws.Range("D56").Value = -0.19
$ws.Range("F56").Value = 959
$ws.Range("G56").Value = 470
$ws.Range("H56").Value = 445

$ws.Range("B57").Value = 9441
$ws.Range("C57").Value = 10844
$ws.Range("D57").Value = -0.13
$ws.Range("E57").Value = 9441
$ws.Range("F57").Value = 10844

$ws.Range("B58").Value = 10932
$ws.Range("C58").Value = 13532
$ws.Range("D58").Value = -0.19
$ws.Range("E58").Value = 10416
$ws.Range("F58").Value = 12897
$ws.Range("G58").Value = 281
$ws.Range("H58").Value = 401

$ws.Range("B59").Value = 21893
$ws.Range("C59").Value = 24293
$ws.Range("D59").Value = -0.099
$ws.Range("E59").Value = 21365
$ws.Range("F59").Value = 23644
$ws.Range("G59").Value = 402
$ws.Range("H59").Value = 488
$ws.Range("K59").Value = 126
$ws.Range("L59").Value = 162

$ws.Range("B60").Value = 3580
$ws.Range("C60").Value = 4163
$ws.Range("D60").Value = -0.14
$ws.Range("F60").Value = 1175
$ws.Range("G60").Value = 2602
$ws.Range("H60").Value = 2919
$ws.Range("K60").Value = 74
$ws.Range("L60").Value = 69

$ws.Range("B61").Value = 66
$ws.Range("C61").Value = 61
$ws.Range("D61").Value = 0.09
$ws.Range("K61").Value = 66
$ws.Range("L61").Value = 61

$ws.Range("C62").Value = 1175
$ws.Range("D62").Value = -0.23
$ws.Range("F62").Value = 1175

$ws.Range("B63").Value = 2609
$ws.Range("C63").Value = 2927
$ws.Range("G63").Value = 2602
$ws.Range("H63").Value = 2919
$ws.Range("K63").Value = 8
$ws.Range("L63").Value = 8

$ws.Range("B64").Value = 1074
$ws.Range("C64").Value = 1084
$ws.Range("D64").Value = -0.009
$ws.Range("E64").Value = 219
$ws.Range("F64").Value = 249
$ws.Range("G64").Value = 806
$ws.Range("H64").Value = 780
$ws.Range("I64").Value = 38
$ws.Range("J64").Value = 40
$ws.Range("K64").Value = 12
$ws.Range("L64").Value = 15

$ws.Range("B65").Value = 398
$ws.Range("C65").Value = 488
$ws.Range("D65").Value = -0.18
$ws.Range("E65").Value = 219
$ws.Range("F65").Value = 249
$ws.Range("G65").Value = 141
$ws.Range("H65").Value = 199
$ws.Range("I65").Value = 38
$ws.Range("J65").Value = 40

$ws.Range("B66").Value = 676
$ws.Range("C66").Value = 596
$ws.Range("D66").Value = 0.13
$ws.Range("G66").Value = 665
$ws.Range("H66").Value = 582
$ws.Range("K66").Value = 12
$ws.Range("L66").Value = 15

$ws.Range("B67").Value = 613093
$ws.Range("C67").Value = 689370
$ws.Range("D67").Value = -0.11
$ws.Range("E67").Value = 450965
$ws.Range("F67").Value = 501628
$ws.Range("G67").Value = 158865
$ws.Range("H67").Value = 183886
$ws.Range("I67").Value = 133
$ws.Range("J67").Value = 149
$ws.Range("K67").Value = 3130
$ws.Range("L67").Value = 3707
